$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new label cells in column D for rows 9, 10, 12
$ws.Range("D9").Value = "Increase in cell size difference that occurred during evolution"
$ws.Range("D10").Value = "Increase in cell size difference due to ftsZ mutation"
$ws.Range("D12").Value = "Proportion of the evolved differences that could potentially be attributed to ftsZ mutation"

# Set column C width to match target (11.7109375 in saved OOXML). This runtime's
# ColumnWidth -> saved-width mapping only lands on multiples of 1/6, so 10.8333
# (-> 11.666666666666666) is the closest achievable approximation to 11.7109375.
$ws.Columns.Item(3).ColumnWidth = 10.8333333333333

# Update selection to C11
$ws.Range("C11").Select()
